$d = $word.ActiveDocument

# Locate the table cell paragraph pair: one paragraph containing exactly
# "Shona" immediately followed by a paragraph containing exactly
# "All add input" (within the same table cell). Word COM Paragraph.Range.Text
# includes trailing control characters (paragraph mark 0x0D, and for the
# last paragraph in a table cell also the cell-mark 0x07), so strip those
# before comparing.
$cr = [char]13
$bell = [char]7

$count = $d.Paragraphs.Count
$shonaIndex = -1
for ($i = 1; $i -lt $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd($cr, $bell)
    if ($text -eq "Shona") {
        $nextPara = $d.Paragraphs.Item($i + 1)
        $nextText = $nextPara.Range.Text.TrimEnd($cr, $bell)
        if ($nextText -eq "All add input") {
            $shonaIndex = $i
            break
        }
    }
}

if ($shonaIndex -eq -1) {
    Write-Output "Target paragraphs not found"
} else {
    # Paragraph 1: "Shona" -> "All" (replace text, keep the paragraph mark
    # and its formatting intact).
    $shonaPara = $d.Paragraphs.Item($shonaIndex)
    $shonaPara.Range.Text = "All"

    # Paragraph 2: remove "All add input" entirely, leaving an empty
    # paragraph (its pPr / paragraph mark stay untouched). Re-fetch the
    # paragraph (positions shifted after the edit above) and build the
    # content-only range from Start + (trimmed text length) rather than
    # subtracting from End, since End snaps past the hidden paragraph /
    # cell-mark run in a way that isn't 1:1 with character counts.
    $inputPara = $d.Paragraphs.Item($shonaIndex + 1)
    $inputRange = $inputPara.Range
    $fullText = $inputRange.Text
    $trimmedText = $fullText.TrimEnd($cr, $bell)
    $textOnlyRange = $d.Range($inputRange.Start, $inputRange.Start + $trimmedText.Length)
    $textOnlyRange.Text = ""

    Write-Output "Edit applied"
}
